# Update the "Avverkningsanmälningar" sheet (Översikt BJUV workbook):
#  - Column C ("Förändrad") bumps from 46063 to 46064 for every data row (2..25).
#  - Rows 7-25 get reordered: for each destination row, the A/B/F/G values
#    (Beteckning / Datum / Markägare / Area (ha)) come from the source row
#    given by the map below (the permutation the upstream scrape produced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow = sourceRow (values read from sourceRow, before any writes)
$rowMap = @{
    7  = 13
    8  = 9
    9  = 11
    10 = 16
    11 = 25
    12 = 14
    13 = 23
    14 = 8
    15 = 17
    16 = 18
    17 = 24
    18 = 19
    19 = 20
    20 = 21
    21 = 22
    22 = 15
    23 = 10
    24 = 12
    25 = 7
}

function Set-CellValue($cell, $val) {
    if ($val -eq "") {
        $cell.Value = ""
    } else {
        $cell.Value = $val
    }
}

# Snapshot the current A/B/F/G values for rows 7..25 before overwriting anything,
# since sources and destinations overlap.
$colA = @{}
$colB = @{}
$colF = @{}
$colG = @{}

for ($r = 7; $r -le 25; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value()
    $colB[$r] = $ws.Cells.Item($r, 2).Value()
    $colF[$r] = $ws.Cells.Item($r, 6).Value()
    $colG[$r] = $ws.Cells.Item($r, 7).Value()
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    Set-CellValue $ws.Cells.Item($destRow, 1) $colA[$srcRow]
    Set-CellValue $ws.Cells.Item($destRow, 2) $colB[$srcRow]
    Set-CellValue $ws.Cells.Item($destRow, 6) $colF[$srcRow]
    Set-CellValue $ws.Cells.Item($destRow, 7) $colG[$srcRow]
}

# Bump column C ("Förändrad") from 46063 to 46064 for every data row (2..25).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 46064
}
